$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldVal = "dnasr281@gmail.com, System"
$newVal = "System, dnasr281@gmail.com"

$ur = $ws.UsedRange
$lastRow = $ur.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $t = $cell.Text
    if ($t -eq $oldVal) {
        $cell.Value = $newVal
    }
}
